{"js": "// Update the worksheet date and all 25 multiplication problems in the table\n// (two-digit x two-digit) from the \"before\" values to the \"after\" values.\nconst replacements = [\n  [\"2023-10-07 Saturday\", \"2023-10-08 Sunday\"],\n  [\"74\u00d724=\", \"69\u00d786=\"],\n  [\"46\u00d737=\", \"23\u00d750=\"],\n  [\"67\u00d717=\", \"23\u00d733=\"],\n  [\"24\u00d746=\", \"91\u00d779=\"],\n  [\"46\u00d715=\", \"82\u00d722=\"],\n  [\"30\u00d729=\", \"96\u00d724=\"],\n  [\"39\u00d773=\", \"31\u00d744=\"],\n  [\"42\u00d762=\", \"72\u00d742=\"],\n  [\"21\u00d798=\", \"17\u00d737=\"],\n  [\"26\u00d773=\", \"54\u00d798=\"],\n  [\"94\u00d780=\", \"65\u00d717=\"],\n  [\"21\u00d770=\", \"18\u00d796=\"],\n  [\"50\u00d785=\", \"63\u00d736=\"],\n  [\"59\u00d780=\", \"88\u00d712=\"],\n  [\"87\u00d765=\", \"14\u00d735=\"],\n  [\"98\u00d781=\", \"95\u00d767=\"],\n  [\"90\u00d798=\", \"37\u00d718=\"],\n  [\"73\u00d712=\", \"89\u00d742=\"],\n  [\"20\u00d729=\", \"31\u00d784=\"],\n  [\"99\u00d765=\", \"34\u00d798=\"],\n  [\"61\u00d721=\", \"49\u00d746=\"],\n  [\"33\u00d764=\", \"95\u00d795=\"],\n  [\"97\u00d755=\", \"17\u00d747=\"],\n  [\"39\u00d793=\", \"49\u00d747=\"],\n  [\"31\u00d716=\", \"79\u00d767=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and all 25 multiplication problems in the table\n# (two-digit x two-digit) from the \"before\" values to the \"after\" values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-10-07 Saturday\", \"2023-10-08 Sunday\"),\n    @(\"74\u00d724=\", \"69\u00d786=\"),\n    @(\"46\u00d737=\", \"23\u00d750=\"),\n    @(\"67\u00d717=\", \"23\u00d733=\"),\n    @(\"24\u00d746=\", \"91\u00d779=\"),\n    @(\"46\u00d715=\", \"82\u00d722=\"),\n    @(\"30\u00d729=\", \"96\u00d724=\"),\n    @(\"39\u00d773=\", \"31\u00d744=\"),\n    @(\"42\u00d762=\", \"72\u00d742=\"),\n    @(\"21\u00d798=\", \"17\u00d737=\"),\n    @(\"26\u00d773=\", \"54\u00d798=\"),\n    @(\"94\u00d780=\", \"65\u00d717=\"),\n    @(\"21\u00d770=\", \"18\u00d796=\"),\n    @(\"50\u00d785=\", \"63\u00d736=\"),\n    @(\"59\u00d780=\", \"88\u00d712=\"),\n    @(\"87\u00d765=\", \"14\u00d735=\"),\n    @(\"98\u00d781=\", \"95\u00d767=\"),\n    @(\"90\u00d798=\", \"37\u00d718=\"),\n    @(\"73\u00d712=\", \"89\u00d742=\"),\n    @(\"20\u00d729=\", \"31\u00d784=\"),\n    @(\"99\u00d765=\", \"34\u00d798=\"),\n    @(\"61\u00d721=\", \"49\u00d746=\"),\n    @(\"33\u00d764=\", \"95\u00d795=\"),\n    @(\"97\u00d755=\", \"17\u00d747=\"),\n    @(\"39\u00d793=\", \"49\u00d747=\"),\n    @(\"31\u00d716=\", \"79\u00d767=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    # wdFindContinue=1, Replace:=wdReplaceAll(2)\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
